$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text-valued numeric columns (Price = D, Volume(1h) = E) ---
# Force Text number format on the specific cells we rewrite so Excel
# does not auto-convert numeric-looking strings (e.g. "1.00", "0.999",
# "0.0000243") into Number cells and silently reformat/round them.
$textCells = @(
    "D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "E7",
    "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "E12", "D13",
    "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18",
    "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23",
    "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28",
    "E28", "D29", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33",
    "E33", "E34", "D35", "E35", "D36", "E36", "D37", "E37", "D38", "E38",
    "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43",
    "D44", "E44", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49",
    "E49", "D50", "E50", "D51", "E51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Cell value updates, row by row ---
# Row 2
$ws.Range("D2").Value = "93.182.61"
$ws.Range("E2").Value = "  +1.72%  "
# Row 3
$ws.Range("D3").Value = "3.095.88"
$ws.Range("E3").Value = "  -0.67%  "
# Row 4
$ws.Range("E4").Value = "  -0.12%  "
# Row 5
$ws.Range("D5").Value = "236.04"
$ws.Range("E5").Value = "  -4.03%  "
# Row 6
$ws.Range("D6").Value = "612.35"
$ws.Range("E6").Value = "  -0.69%  "
# Row 7
$ws.Range("E7").Value = "  +3.08%  "
# Row 8
$ws.Range("D8").Value = "0.387"
$ws.Range("E8").Value = "  +1.25%  "
# Row 9
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.10%  "
# Row 10
$ws.Range("D10").Value = "0.827"
$ws.Range("E10").Value = "  +13.13%  "
# Row 11
$ws.Range("D11").Value = "3.096.11"
$ws.Range("E11").Value = "  -0.61%  "
# Row 12
$ws.Range("E12").Value = "  -2.83%  "
# Row 13
$ws.Range("D13").Value = "0.0000243"
$ws.Range("E13").Value = "  -2.86%  "
# Row 14
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "34.84"
$ws.Range("E14").Value = "  +0.34%  "
# Row 15
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "93.008.41"
$ws.Range("E15").Value = "  +1.58%  "
# Row 16
$ws.Range("D16").Value = "5.40"
$ws.Range("E16").Value = "  -3.30%  "
# Row 17
$ws.Range("D17").Value = "3.671.55"
$ws.Range("E17").Value = "  -0.83%  "
# Row 18
$ws.Range("D18").Value = "3.119.12"
$ws.Range("E18").Value = "  +0.96%  "
# Row 19
$ws.Range("D19").Value = "3.67"
$ws.Range("E19").Value = "  -0.61%  "
# Row 20
$ws.Range("D20").Value = "14.62"
$ws.Range("E20").Value = "  -1.28%  "
# Row 21
$ws.Range("D21").Value = "5.94"
$ws.Range("E21").Value = "  +2.54%  "
# Row 22
$ws.Range("D22").Value = "441.00"
$ws.Range("E22").Value = "  -1.18%  "
# Row 23
$ws.Range("D23").Value = "0.0000198"
$ws.Range("E23").Value = "  -1.44%  "
# Row 24
$ws.Range("D24").Value = "9.04"
$ws.Range("E24").Value = "  -4.49%  "
# Row 25
$ws.Range("D25").Value = "8.22"
$ws.Range("E25").Value = "  +4.34%  "
# Row 26
$ws.Range("D26").Value = "5.67"
$ws.Range("E26").Value = "  -2.82%  "
# Row 27
$ws.Range("D27").Value = "12.65"
$ws.Range("E27").Value = "  +7.97%  "
# Row 28
$ws.Range("D28").Value = "85.62"
$ws.Range("E28").Value = "  -2.54%  "
# Row 29
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.24%  "
# Row 30
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "0.250"
$ws.Range("E30").Value = "  +6.81%  "
# Row 31
$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").Value = "0.181"
$ws.Range("E31").Value = "  +8.40%  "
# Row 32
$ws.Range("D32").Value = "0.124"
$ws.Range("E32").Value = "  -14.14%  "
# Row 33
$ws.Range("D33").Value = "9.18"
$ws.Range("E33").Value = "  -1.14%  "
# Row 34
$ws.Range("E34").Value = "  +0.69%  "
# Row 35
$ws.Range("D35").Value = "7.90"
$ws.Range("E35").Value = "  +0.73%  "
# Row 36
$ws.Range("D36").Value = "0.159"
$ws.Range("E36").Value = "  -10.06%  "
# Row 37
$ws.Range("D37").Value = "25.89"
$ws.Range("E37").Value = "  -1.15%  "
# Row 38
$ws.Range("D38").Value = "3.99"
$ws.Range("E38").Value = "  -4.66%  "
# Row 39
$ws.Range("D39").Value = "1.90"
$ws.Range("E39").Value = "  -2.11%  "
# Row 40
$ws.Range("D40").Value = "23.97"
$ws.Range("E40").Value = "  +8.08%  "
# Row 41
$ws.Range("D41").Value = "0.441"
$ws.Range("E41").Value = "  +0.66%  "
# Row 42
$ws.Range("D42").Value = "1.29"
$ws.Range("E42").Value = "  -1.18%  "
# Row 43
$ws.Range("D43").Value = "473.87"
$ws.Range("E43").Value = "  -3.25%  "
# Row 44
$ws.Range("D44").Value = "3.26"
$ws.Range("E44").Value = "  -4.06%  "
# Row 45
$ws.Range("E45").Value = "  +0.03%  "
# Row 46
$ws.Range("D46").Value = "159.01"
$ws.Range("E46").Value = "  +0.84%  "
# Row 47
$ws.Range("D47").Value = "0.697"
$ws.Range("E47").Value = "  -1.10%  "
# Row 48
$ws.Range("D48").Value = "1.86"
$ws.Range("E48").Value = "  -2.53%  "
# Row 49
$ws.Range("D49").Value = "1.32"
$ws.Range("E49").Value = "  -2.26%  "
# Row 50
$ws.Range("D50").Value = "43.86"
$ws.Range("E50").Value = "  -0.44%  "
# Row 51
$ws.Range("D51").Value = "4.36"
$ws.Range("E51").Value = "  -0.40%  "
